$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B20").Value = 18
$ws.Range("D20").Value = "Mesh"
$ws.Range("E20").Value = "Load"
$ws.Range("C20").Value = "Error Loading Mesh"

$ws.Range("B21").Value = 19
$ws.Range("C21").Value = "Pool Allocation Service Unknown Error"
$ws.Range("D21").Value = "PoolAllocationService"
$ws.Range("E21").Value = "MakeUniquePoolPtr"

$ws.Range("C26").Select()
